$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("A6").Value = 130746493
$ws.Range("B6").Value = 57881
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 100049
$ws.Range("F6").Value = "Spillkråka"
$ws.Range("G6").Value = "Dryocopus martius"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("M6").Value = "äldre spår"
$ws.Range("Q6").Value = 447766
$ws.Range("R6").Value = 6784433

# Row 8
$ws.Range("A8").Value = 130746520
$ws.Range("B8").Value = 8451
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 106545
$ws.Range("F8").Value = "Mindre märgborre"
$ws.Range("G8").Value = "Tomicus minor"
$ws.Range("H8").Value = "(Hartig, 1834)"
$ws.Range("M8").Value = "äldre gnagspår"
$ws.Range("Q8").Value = 447846
$ws.Range("R8").Value = 6784643

# Row 10
$ws.Range("A10").Value = 130746562
$ws.Range("B10").Value = 79243
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("M10").Value = ""
$ws.Range("Q10").Value = 447730
$ws.Range("R10").Value = 6784717

# Row 11
$ws.Range("A11").Value = 130746524
$ws.Range("Q11").Value = 447932
$ws.Range("R11").Value = 6784551

# Row 12
$ws.Range("A12").Value = 130746530
$ws.Range("B12").Value = 8451
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 106545
$ws.Range("F12").Value = "Mindre märgborre"
$ws.Range("G12").Value = "Tomicus minor"
$ws.Range("H12").Value = "(Hartig, 1834)"
$ws.Range("M12").Value = "äldre gnagspår"
$ws.Range("Q12").Value = 447855
$ws.Range("R12").Value = 6784599

# Row 13
$ws.Range("A13").Value = 130746519
$ws.Range("M13").Value = "färska gnagspår"
$ws.Range("Q13").Value = 447826
$ws.Range("R13").Value = 6784623

# Row 16
$ws.Range("A16").Value = 130746554
$ws.Range("B16").Value = 79243
$ws.Range("E16").Value = 6425
$ws.Range("F16").Value = "Garnlav"
$ws.Range("G16").Value = "Alectoria sarmentosa"
$ws.Range("H16").Value = "(Ach.) Ach."
$ws.Range("M16").Value = ""
$ws.Range("Q16").Value = 447915
$ws.Range("R16").Value = 6784490

# Row 17
$ws.Range("A17").Value = 130746525
$ws.Range("B17").Value = 8451
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 106545
$ws.Range("F17").Value = "Mindre märgborre"
$ws.Range("G17").Value = "Tomicus minor"
$ws.Range("H17").Value = "(Hartig, 1834)"
$ws.Range("M17").Value = "äldre gnagspår"
$ws.Range("Q17").Value = 447933
$ws.Range("R17").Value = 6784552

# Row 18
$ws.Range("A18").Value = 130746499
$ws.Range("Q18").Value = 447930
$ws.Range("R18").Value = 6784568

# Row 19
$ws.Range("A19").Value = 130746500
$ws.Range("B19").Value = 57881
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 100049
$ws.Range("F19").Value = "Spillkråka"
$ws.Range("G19").Value = "Dryocopus martius"
$ws.Range("H19").Value = "(Linnaeus, 1758)"
$ws.Range("M19").Value = "äldre spår"
$ws.Range("Q19").Value = 447888
$ws.Range("R19").Value = 6784627

# Row 20
$ws.Range("A20").Value = 130746546
$ws.Range("B20").Value = 92246
$ws.Range("E20").Value = 5420
$ws.Range("F20").Value = "Grovticka"
$ws.Range("G20").Value = "Phaeolus schweinitzii"
$ws.Range("H20").Value = "(Fr.) Pat."
$ws.Range("M20").Value = ""
$ws.Range("Q20").Value = 447898
$ws.Range("R20").Value = 6784612

# Row 21
$ws.Range("A21").Value = 130746490
$ws.Range("B21").Value = 58043
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 103021
$ws.Range("F21").Value = "Talltita"
$ws.Range("G21").Value = "Poecile montanus"
$ws.Range("H21").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("M21").Value = "lockläte, övriga läten"
$ws.Range("Q21").Value = 447888
$ws.Range("R21").Value = 6784627

# Row 22
$ws.Range("A22").Value = 130746531
$ws.Range("B22").Value = 8451
$ws.Range("D22").Value = "LC"
$ws.Range("E22").Value = 106545
$ws.Range("F22").Value = "Mindre märgborre"
$ws.Range("G22").Value = "Tomicus minor"
$ws.Range("H22").Value = "(Hartig, 1834)"
$ws.Range("M22").Value = "äldre gnagspår"
$ws.Range("Q22").Value = 447932
$ws.Range("R22").Value = 6784555

# Row 23
$ws.Range("A23").Value = 130746556
$ws.Range("B23").Value = 79243
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 6425
$ws.Range("F23").Value = "Garnlav"
$ws.Range("G23").Value = "Alectoria sarmentosa"
$ws.Range("H23").Value = "(Ach.) Ach."
$ws.Range("M23").Value = ""
$ws.Range("Q23").Value = 447876
$ws.Range("R23").Value = 6784529

# Row 24
$ws.Range("A24").Value = 130746510
$ws.Range("B24").Value = 8451
$ws.Range("D24").Value = "LC"
$ws.Range("E24").Value = 106545
$ws.Range("F24").Value = "Mindre märgborre"
$ws.Range("G24").Value = "Tomicus minor"
$ws.Range("H24").Value = "(Hartig, 1834)"
$ws.Range("M24").Value = "äldre gnagspår"
$ws.Range("Q24").Value = 447718
$ws.Range("R24").Value = 6784441

# Row 25
$ws.Range("A25").Value = 130746564
$ws.Range("B25").Value = 79243
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 6425
$ws.Range("F25").Value = "Garnlav"
$ws.Range("G25").Value = "Alectoria sarmentosa"
$ws.Range("H25").Value = "(Ach.) Ach."
$ws.Range("M25").Value = ""
$ws.Range("Q25").Value = 447866
$ws.Range("R25").Value = 6784648

# Row 26
$ws.Range("A26").Value = 130746569
$ws.Range("B26").Value = 79243
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 6425
$ws.Range("F26").Value = "Garnlav"
$ws.Range("G26").Value = "Alectoria sarmentosa"
$ws.Range("H26").Value = "(Ach.) Ach."
$ws.Range("M26").Value = ""
$ws.Range("Q26").Value = 447856
$ws.Range("R26").Value = 6784518

# Row 27
$ws.Range("A27").Value = 130746515
$ws.Range("B27").Value = 8451
$ws.Range("D27").Value = "LC"
$ws.Range("E27").Value = 106545
$ws.Range("F27").Value = "Mindre märgborre"
$ws.Range("G27").Value = "Tomicus minor"
$ws.Range("H27").Value = "(Hartig, 1834)"
$ws.Range("M27").Value = "färska gnagspår"
$ws.Range("Q27").Value = 447716
$ws.Range("R27").Value = 6784496

# Row 28
$ws.Range("A28").Value = 130746518
$ws.Range("B28").Value = 8451
$ws.Range("D28").Value = "LC"
$ws.Range("E28").Value = 106545
$ws.Range("F28").Value = "Mindre märgborre"
$ws.Range("G28").Value = "Tomicus minor"
$ws.Range("H28").Value = "(Hartig, 1834)"
$ws.Range("M28").Value = "äldre gnagspår"
$ws.Range("Q28").Value = 447815
$ws.Range("R28").Value = 6784612

# Row 32
$ws.Range("A32").Value = 130746523
$ws.Range("B32").Value = 8451
$ws.Range("D32").Value = "LC"
$ws.Range("E32").Value = 106545
$ws.Range("F32").Value = "Mindre märgborre"
$ws.Range("G32").Value = "Tomicus minor"
$ws.Range("H32").Value = "(Hartig, 1834)"
$ws.Range("M32").Value = "äldre gnagspår"
$ws.Range("Q32").Value = 447940
$ws.Range("R32").Value = 6784589

# Row 34
$ws.Range("A34").Value = 130746555
$ws.Range("B34").Value = 79243
$ws.Range("D34").Value = "NT"
$ws.Range("E34").Value = 6425
$ws.Range("F34").Value = "Garnlav"
$ws.Range("G34").Value = "Alectoria sarmentosa"
$ws.Range("H34").Value = "(Ach.) Ach."
$ws.Range("M34").Value = ""
$ws.Range("Q34").Value = 447906
$ws.Range("R34").Value = 6784505

# Row 35
$ws.Range("A35").Value = 130746497
$ws.Range("B35").Value = 57881
$ws.Range("E35").Value = 100049
$ws.Range("F35").Value = "Spillkråka"
$ws.Range("G35").Value = "Dryocopus martius"
$ws.Range("H35").Value = "(Linnaeus, 1758)"
$ws.Range("M35").Value = "färska spår"
$ws.Range("Q35").Value = 447838
$ws.Range("R35").Value = 6784644

# Row 36
$ws.Range("A36").Value = 130746560
$ws.Range("Q36").Value = 447685
$ws.Range("R36").Value = 6784529

# Row 37
$ws.Range("A37").Value = 130746565
$ws.Range("Q37").Value = 447912
$ws.Range("R37").Value = 6784599

# Row 41
$ws.Range("A41").Value = 130746563
$ws.Range("B41").Value = 79243
$ws.Range("E41").Value = 6425
$ws.Range("F41").Value = "Garnlav"
$ws.Range("G41").Value = "Alectoria sarmentosa"
$ws.Range("H41").Value = "(Ach.) Ach."
$ws.Range("M41").Value = ""
$ws.Range("Q41").Value = 447832
$ws.Range("R41").Value = 6784636

# Row 42
$ws.Range("A42").Value = 130746566
$ws.Range("Q42").Value = 447949
$ws.Range("R42").Value = 6784550

# Row 43
$ws.Range("A43").Value = 130746495
$ws.Range("B43").Value = 57881
$ws.Range("D43").Value = "NT"
$ws.Range("E43").Value = 100049
$ws.Range("F43").Value = "Spillkråka"
$ws.Range("G43").Value = "Dryocopus martius"
$ws.Range("H43").Value = "(Linnaeus, 1758)"
$ws.Range("M43").Value = "färska spår"
$ws.Range("Q43").Value = 447746
$ws.Range("R43").Value = 6784474

# Row 44
$ws.Range("A44").Value = 130746511
$ws.Range("B44").Value = 8451
$ws.Range("D44").Value = "LC"
$ws.Range("E44").Value = 106545
$ws.Range("F44").Value = "Mindre märgborre"
$ws.Range("G44").Value = "Tomicus minor"
$ws.Range("H44").Value = "(Hartig, 1834)"
$ws.Range("M44").Value = "äldre gnagspår"
$ws.Range("Q44").Value = 447748
$ws.Range("R44").Value = 6784472
